$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "Associated Error"

$ws.Range("E2").Formula = "=A3"
$ws.Range("E3").Formula = "=A6"
$ws.Range("E4").Formula = "=A9"
$ws.Range("E5").Formula = "=A12"
$ws.Range("E6").Formula = "=A15"
$ws.Range("E7").Formula = "=A18"
$ws.Range("E8").Formula = "=A21"
$ws.Range("E9").Formula = "=A24"
$ws.Range("E10").Formula = "=A27"
$ws.Range("E11").Formula = "=A30"
$ws.Range("E12").Formula = "=A33"
$ws.Range("E13").Formula = "=A36"

$ws.Columns("D").ColumnWidth = 10.833333333333332
$ws.Columns("E").ColumnWidth = 14.666666666666668

$ws.Range("E14").Select()
